# CarScenarios.xlsx edit script
#
# Summary of the intended change (per commit message / diff):
#  - "fix for AddToCart method and added BookNow function for car product"
#  - The ExecutionPipeline column (column B, row 2) on every scenario sheet
#    is extended from "Search|AddToCart" to
#    "Search|AddToCart|checkoutTrip|Login|enterpaxinfo|confirmpaxinfo".
#  - Two of the sheets (Car_Airport_To_City and Car_City_To_Airport) go one
#    step further and use the "BookNow"/pay flow variant:
#    "Search|AddToCart|checkoutTrip|Login|enterpaxinfo|confirmpaxinfo||Paynow".
#  - The previously active sheet/tab (Car_City_To_SameAsPickUp) is no longer
#    the active tab; Car_City_To_Airport becomes the active tab instead.
#  - Assorted cursor/selection bookkeeping left behind by whoever was
#    clicking around the workbook while making the edits.

$wb = $excel.ActiveWorkbook

# Common new ExecutionPipeline text used by every sheet.
$basePipeline = "Search|AddToCart|checkoutTrip|Login|enterpaxinfo|confirmpaxinfo"
# Extended variant (adds the Paynow/BookNow step) used by two sheets.
$payNowPipeline = "Search|AddToCart|checkoutTrip|Login|enterpaxinfo|confirmpaxinfo||Paynow"

$wsAirportToAirport   = $wb.Worksheets.Item("Car_Airport_To_Airport")
$wsAirportToCity      = $wb.Worksheets.Item("Car_Airport_To_City")
$wsAirportToSameAsPU  = $wb.Worksheets.Item("Car_Airport_To_SameAsPickUp")
$wsCityToAirport      = $wb.Worksheets.Item("Car_City_To_Airport")
$wsCityToCity         = $wb.Worksheets.Item("Car_City_To_City")
$wsCityToSameAsPU     = $wb.Worksheets.Item("Car_City_To_SameAsPickUp")

# Update the ExecutionPipeline (column B, row 2) value on every sheet.
$wsAirportToAirport.Range("B2").Value  = $basePipeline
$wsAirportToCity.Range("B2").Value     = $payNowPipeline
$wsAirportToSameAsPU.Range("B2").Value = $basePipeline
$wsCityToAirport.Range("B2").Value     = $payNowPipeline
$wsCityToCity.Range("B2").Value        = $basePipeline
$wsCityToSameAsPU.Range("B2").Value    = $basePipeline

# Leftover cell-cursor positions on the non-active sheets.
$wsCityToSameAsPU.Range("B2").Select()

$wsAirportToAirport.Activate()
$wsAirportToAirport.Range("B7").Select()

$wsAirportToCity.Activate()
$wsAirportToCity.Range("B11").Select()

$wsCityToCity.Activate()
$wsCityToCity.Range("B13").Select()

# Car_City_To_Airport becomes the active tab/sheet, with its cursor on C8.
$wsCityToAirport.Activate()
$wsCityToAirport.Range("C8").Select()
